# Config.xlsx edit - "Made changes in config file using assets"
#
# Adds three Orchestrator-Asset rows to the Assets sheet (mirroring the
# Workbook1_Path / Workbook2_Path / Workbook3_Path entries already present
# on the Settings sheet), widens the two leftmost columns on Assets to fit
# the new text, sets the Assets sheet's print orientation to portrait, and
# leaves the Assets sheet as the active tab/selection (it was the sheet the
# user finished working on), while Settings keeps its own last selection.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$assets   = $wb.Worksheets.Item("Assets")

# --- Settings sheet: this is where the user started, reviewing/selecting
#     B9 (Workbook3_Path's value) before moving over to the Assets sheet.
$settings.Activate()
$settings.Range("B9").Select()

# --- Assets sheet: add the three new Forward_* asset rows under the header.
$assets.Activate()

$assets.Range("A2").Value = "CompCodeTracker_FilePath"
$assets.Range("B2").Value = "Forward_CompCodeTracker_FilePath"
$assets.Range("C2").Value = "Feat_Forward"

$assets.Range("A3").Value = "COMP Passes Master_FilePath"
$assets.Range("B3").Value = "Forward_COMP Passes Master_FilePath"
$assets.Range("C3").Value = "Feat_Forward"

$assets.Range("A4").Value = "FUSIONRegistrationReport_FilePath"
$assets.Range("B4").Value = "Forward_FUSIONRegistrationReport_FilePath"
$assets.Range("C4").Value = "Feat_Forward"

# Widen column A (best-fit to the new, longer asset names) and column B
# (wide enough for the longest "Forward_..." label).
$assets.Columns.Item(1).AutoFit()
$assets.Columns.Item(2).ColumnWidth = 49

# Print setup: portrait orientation for the Assets sheet.
$assets.PageSetup.Orientation = 1

# Final selection / active sheet: Assets, cell A2.
$assets.Range("A2").Select()
